$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Cells.Item(1,1).Value = "venue"
$ws.Cells.Item(1,2).Value = "date"
$ws.Cells.Item(1,3).Value = "result"
$ws.Cells.Item(1,4).Value = "ownTeam"
$ws.Cells.Item(1,5).Value = "oppTeam"
$ws.Cells.Item(1,6).Value = "batsman"
$ws.Cells.Item(1,7).Value = "totalRuns"
$ws.Cells.Item(1,8).Value = "totalBalls"
$ws.Cells.Item(1,9).Value = "total4s"
$ws.Cells.Item(1,10).Value = "total6s"
$ws.Cells.Item(1,11).Value = "sr"

# Row 2
$ws.Cells.Item(2,1).Value = " Abu Dhabi"
$ws.Cells.Item(2,2).Value = " October 25 2020"
$ws.Cells.Item(2,3).Value = "Royals won by 8 wickets (with 10 balls remaining)"
$ws.Cells.Item(2,4).Value = "Rajasthan Royals"
$ws.Cells.Item(2,5).Value = "Mumbai Indians"
$ws.Cells.Item(2,6).Value = "Ben Stokes "
$ws.Cells.Item(2,7).Value = "'107"
$ws.Cells.Item(2,8).Value = "'60"
$ws.Cells.Item(2,9).Value = "'14"
$ws.Cells.Item(2,10).Value = "'3"
$ws.Cells.Item(2,11).Value = "'178.33"

# Row 3
$ws.Cells.Item(3,1).Value = " Dubai (DSC)"
$ws.Cells.Item(3,2).Value = " October 17 2020"
$ws.Cells.Item(3,3).Value = "RCB won by 7 wickets (with 2 balls remaining)"
$ws.Cells.Item(3,4).Value = "Rajasthan Royals"
$ws.Cells.Item(3,5).Value = "Royal Challengers Bangalore"
$ws.Cells.Item(3,6).Value = "Ben Stokes "
$ws.Cells.Item(3,7).Value = "'15"
$ws.Cells.Item(3,8).Value = "'19"
$ws.Cells.Item(3,9).Value = "'2"
$ws.Cells.Item(3,10).Value = "'0"
$ws.Cells.Item(3,11).Value = "'78.94"

# Row 4
$ws.Cells.Item(4,1).Value = " Abu Dhabi"
$ws.Cells.Item(4,2).Value = " October 30 2020"
$ws.Cells.Item(4,3).Value = "Royals won by 7 wickets (with 15 balls remaining)"
$ws.Cells.Item(4,4).Value = "Rajasthan Royals"
$ws.Cells.Item(4,5).Value = "Kings XI Punjab"
$ws.Cells.Item(4,6).Value = "Ben Stokes "
$ws.Cells.Item(4,7).Value = "'50"
$ws.Cells.Item(4,8).Value = "'26"
$ws.Cells.Item(4,9).Value = "'6"
$ws.Cells.Item(4,10).Value = "'3"
$ws.Cells.Item(4,11).Value = "'192.30"

# Row 5
$ws.Cells.Item(5,1).Value = " Dubai (DSC)"
$ws.Cells.Item(5,2).Value = " October 14 2020"
$ws.Cells.Item(5,3).Value = "Capitals won by 13 runs"
$ws.Cells.Item(5,4).Value = "Rajasthan Royals"
$ws.Cells.Item(5,5).Value = "Delhi Capitals"
$ws.Cells.Item(5,6).Value = "Ben Stokes "
$ws.Cells.Item(5,7).Value = "'41"
$ws.Cells.Item(5,8).Value = "'35"
$ws.Cells.Item(5,9).Value = "'6"
$ws.Cells.Item(5,10).Value = "'0"
$ws.Cells.Item(5,11).Value = "'117.14"

# Row 6
$ws.Cells.Item(6,1).Value = " Dubai (DSC)"
$ws.Cells.Item(6,2).Value = " October 22 2020"
$ws.Cells.Item(6,3).Value = "Sunrisers won by 8 wickets (with 11 balls remaining)"
$ws.Cells.Item(6,4).Value = "Rajasthan Royals"
$ws.Cells.Item(6,5).Value = "Sunrisers Hyderabad"
$ws.Cells.Item(6,6).Value = "Ben Stokes "
$ws.Cells.Item(6,7).Value = "'30"
$ws.Cells.Item(6,8).Value = "'32"
$ws.Cells.Item(6,9).Value = "'2"
$ws.Cells.Item(6,10).Value = "'0"
$ws.Cells.Item(6,11).Value = "'93.75"

# Row 7
$ws.Cells.Item(7,1).Value = " Dubai (DSC)"
$ws.Cells.Item(7,2).Value = " November 01 2020"
$ws.Cells.Item(7,3).Value = "KKR won by 60 runs"
$ws.Cells.Item(7,4).Value = "Rajasthan Royals"
$ws.Cells.Item(7,5).Value = "Kolkata Knight Riders"
$ws.Cells.Item(7,6).Value = "Ben Stokes "
$ws.Cells.Item(7,7).Value = "'18"
$ws.Cells.Item(7,8).Value = "'11"
$ws.Cells.Item(7,9).Value = "'2"
$ws.Cells.Item(7,10).Value = "'1"
$ws.Cells.Item(7,11).Value = "'163.63"

# Row 8
$ws.Cells.Item(8,1).Value = " Abu Dhabi"
$ws.Cells.Item(8,2).Value = " October 19 2020"
$ws.Cells.Item(8,3).Value = "Royals won by 7 wickets (with 15 balls remaining)"
$ws.Cells.Item(8,4).Value = "Rajasthan Royals"
$ws.Cells.Item(8,5).Value = "Chennai Super Kings"
$ws.Cells.Item(8,6).Value = "Ben Stokes "
$ws.Cells.Item(8,7).Value = "'19"
$ws.Cells.Item(8,8).Value = "'11"
$ws.Cells.Item(8,9).Value = "'3"
$ws.Cells.Item(8,10).Value = "'0"
$ws.Cells.Item(8,11).Value = "'172.72"

# Row 9
$ws.Cells.Item(9,1).Value = " Dubai (DSC)"
$ws.Cells.Item(9,2).Value = " October 11 2020"
$ws.Cells.Item(9,3).Value = "Royals won by 5 wickets (with 1 ball remaining)"
$ws.Cells.Item(9,4).Value = "Rajasthan Royals"
$ws.Cells.Item(9,5).Value = "Sunrisers Hyderabad"
$ws.Cells.Item(9,6).Value = "Ben Stokes "
$ws.Cells.Item(9,7).Value = "'5"
$ws.Cells.Item(9,8).Value = "'6"
$ws.Cells.Item(9,9).Value = "'1"
$ws.Cells.Item(9,10).Value = "'0"
$ws.Cells.Item(9,11).Value = "'83.33"
